# Update existing rows 2-7 and append new rows 8-10 per the NATMI re-run
# ("Natmi following Dr Hou advice") for the Sfrp1-Fzd2 ligand-receptor pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sfrp1"
$ws.Cells.Item(2,3).Value = "Fzd2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.3669776666666666
$ws.Cells.Item(2,8).Value = 1.100933
$ws.Cells.Item(2,9).Value = 0.0181959334720815
$ws.Cells.Item(2,10).Value = 0.0181959334720815
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.07629866666666667
$ws.Cells.Item(2,14).Value = 0.228896
$ws.Cells.Item(2,15).Value = 0.004108848954870246
$ws.Cells.Item(2,16).Value = 0.004108848954870246
$ws.Cells.Item(2,17).Value = 0.02799990666311111
$ws.Cells.Item(2,18).Value = 0.251999159968
$ws.Cells.Item(2,19).Value = 0.00007476434222965059
$ws.Cells.Item(2,20).Value = 0.00007476434222965059

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sfrp1"
$ws.Cells.Item(3,3).Value = "Fzd2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.3669776666666666
$ws.Cells.Item(3,8).Value = 1.100933
$ws.Cells.Item(3,9).Value = 0.0181959334720815
$ws.Cells.Item(3,10).Value = 0.0181959334720815
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 12.99468133333333
$ws.Cells.Item(3,14).Value = 38.984044
$ws.Cells.Item(3,15).Value = 0.6997918200668237
$ws.Cells.Item(3,16).Value = 0.6997918200668237
$ws.Cells.Item(3,17).Value = 4.768757834783555
$ws.Cells.Item(3,18).Value = 42.918820513052
$ws.Cells.Item(3,19).Value = 0.01273336540224275
$ws.Cells.Item(3,20).Value = 0.01273336540224275

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sfrp1"
$ws.Cells.Item(4,3).Value = "Fzd2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.3669776666666666
$ws.Cells.Item(4,8).Value = 1.100933
$ws.Cells.Item(4,9).Value = 0.0181959334720815
$ws.Cells.Item(4,10).Value = 0.0181959334720815
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.498373000000001
$ws.Cells.Item(4,14).Value = 16.495119
$ws.Cells.Item(4,15).Value = 0.2960993309783061
$ws.Cells.Item(4,16).Value = 0.2960993309783061
$ws.Cells.Item(4,17).Value = 2.017780094003
$ws.Cells.Item(4,18).Value = 18.160020846027
$ws.Cells.Item(4,19).Value = 0.005387803727609097
$ws.Cells.Item(4,20).Value = 0.005387803727609097

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Sfrp1"
$ws.Cells.Item(5,3).Value = "Fzd2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 17.23456
$ws.Cells.Item(5,8).Value = 51.70368000000001
$ws.Cells.Item(5,9).Value = 0.8545449373774706
$ws.Cells.Item(5,10).Value = 0.8545449373774706
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.07629866666666667
$ws.Cells.Item(5,14).Value = 0.228896
$ws.Cells.Item(5,15).Value = 0.004108848954870246
$ws.Cells.Item(5,16).Value = 0.004108848954870246
$ws.Cells.Item(5,17).Value = 1.314973948586667
$ws.Cells.Item(5,18).Value = 11.83476553728
$ws.Cells.Item(5,19).Value = 0.00351119607283308
$ws.Cells.Item(5,20).Value = 0.00351119607283308

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sfrp1"
$ws.Cells.Item(6,3).Value = "Fzd2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 17.23456
$ws.Cells.Item(6,8).Value = 51.70368000000001
$ws.Cells.Item(6,9).Value = 0.8545449373774706
$ws.Cells.Item(6,10).Value = 0.8545449373774706
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 12.99468133333333
$ws.Cells.Item(6,14).Value = 38.984044
$ws.Cells.Item(6,15).Value = 0.6997918200668237
$ws.Cells.Item(6,16).Value = 0.6997918200668237
$ws.Cells.Item(6,17).Value = 223.9576151202133
$ws.Cells.Item(6,18).Value = 2015.61853608192
$ws.Cells.Item(6,19).Value = 0.5980035570562701
$ws.Cells.Item(6,20).Value = 0.5980035570562701

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sfrp1"
$ws.Cells.Item(7,3).Value = "Fzd2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 17.23456
$ws.Cells.Item(7,8).Value = 51.70368000000001
$ws.Cells.Item(7,9).Value = 0.8545449373774706
$ws.Cells.Item(7,10).Value = 0.8545449373774706
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.498373000000001
$ws.Cells.Item(7,14).Value = 16.495119
$ws.Cells.Item(7,15).Value = 0.2960993309783061
$ws.Cells.Item(7,16).Value = 0.2960993309783061
$ws.Cells.Item(7,17).Value = 94.76203937088002
$ws.Cells.Item(7,18).Value = 852.8583543379202
$ws.Cells.Item(7,19).Value = 0.2530301842483675
$ws.Cells.Item(7,20).Value = 0.2530301842483675

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Sfrp1"
$ws.Cells.Item(8,3).Value = "Fzd2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.566576666666667
$ws.Cells.Item(8,8).Value = 7.69973
$ws.Cells.Item(8,9).Value = 0.1272591291504479
$ws.Cells.Item(8,10).Value = 0.1272591291504479
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.07629866666666667
$ws.Cells.Item(8,14).Value = 0.228896
$ws.Cells.Item(8,15).Value = 0.004108848954870246
$ws.Cells.Item(8,16).Value = 0.004108848954870246
$ws.Cells.Item(8,17).Value = 0.1958263775644445
$ws.Cells.Item(8,18).Value = 1.76243739808
$ws.Cells.Item(8,19).Value = 0.0005228885398075156
$ws.Cells.Item(8,20).Value = 0.0005228885398075155

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Sfrp1"
$ws.Cells.Item(9,3).Value = "Fzd2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.566576666666667
$ws.Cells.Item(9,8).Value = 7.69973
$ws.Cells.Item(9,9).Value = 0.1272591291504479
$ws.Cells.Item(9,10).Value = 0.1272591291504479
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 12.99468133333333
$ws.Cells.Item(9,14).Value = 38.984044
$ws.Cells.Item(9,15).Value = 0.6997918200668237
$ws.Cells.Item(9,16).Value = 0.6997918200668237
$ws.Cells.Item(9,17).Value = 33.35184590090222
$ws.Cells.Item(9,18).Value = 300.16661310812
$ws.Cells.Item(9,19).Value = 0.08905489760831094
$ws.Cells.Item(9,20).Value = 0.08905489760831092

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sfrp1"
$ws.Cells.Item(10,3).Value = "Fzd2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.566576666666667
$ws.Cells.Item(10,8).Value = 7.69973
$ws.Cells.Item(10,9).Value = 0.1272591291504479
$ws.Cells.Item(10,10).Value = 0.1272591291504479
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.498373000000001
$ws.Cells.Item(10,14).Value = 16.495119
$ws.Cells.Item(10,15).Value = 0.2960993309783061
$ws.Cells.Item(10,16).Value = 0.2960993309783061
$ws.Cells.Item(10,17).Value = 14.11199584643
$ws.Cells.Item(10,18).Value = 127.00796261787
$ws.Cells.Item(10,19).Value = 0.03768134300232948
$ws.Cells.Item(10,20).Value = 0.03768134300232948
